$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, matching the style used by the other headers (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Values for the new "Save" column (rows 2-9)
$saveValues = @(0, 0, 0, 0, 1, 0, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
